$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Range("H28").Value = 425.84616
$ws.Range("J28").Value = 615.6667
$ws.Range("L28").Value = 615.6667
$ws.Range("N28").Value = -1585.6667
# Row 33
$ws.Range("H33").Value = 330.64285
$ws.Range("I33").Value = 371.125
$ws.Range("J33").Value = 276.66666
$ws.Range("K33").Value = 371.125
$ws.Range("L33").Value = 276.66666
$ws.Range("M33").Value = -142.125
$ws.Range("N33").Value = -734.66666
# Row 86
$ws.Range("H86").Value = 2750.2222
$ws.Range("I86").Value = 2643.2856
$ws.Range("K86").Value = 2643.2856
$ws.Range("M86").Value = -1520.2856
# Row 89
$ws.Range("H89").Value = 2750.2222
$ws.Range("I89").Value = 2643.2856
$ws.Range("K89").Value = 13216.428
$ws.Range("M89").Value = -7600.428
# Row 93
$ws.Range("H93").Value = 46000
$ws.Range("J93").Value = 46000
$ws.Range("L93").Value = 46000
$ws.Range("N93").Value = -50992
# Row 96
$ws.Range("H96").Value = 774.6
$ws.Range("I96").Value = 449
$ws.Range("J96").Value = 991.6667
$ws.Range("K96").Value = 1347
$ws.Range("L96").Value = 2975.0001
$ws.Range("M96").Value = 26
$ws.Range("N96").Value = -5721.0001
# Row 98
$ws.Range("H98").Value = 2375.1365
$ws.Range("I98").Value = 2119.1052
$ws.Range("K98").Value = 2119.1052
$ws.Range("M98").Value = -621.1052
# Row 103
$ws.Range("H103").Value = 1387.5834
$ws.Range("I103").Value = 1364.5
$ws.Range("J103").Value = 1433.75
$ws.Range("K103").Value = 4093.5
$ws.Range("L103").Value = 4301.25
$ws.Range("M103").Value = -3507.5
$ws.Range("N103").Value = -5473.25
# Row 106
$ws.Range("H106").Value = 2077.1428
$ws.Range("I106").Value = 1923.3334
$ws.Range("K106").Value = 1923.3334
$ws.Range("M106").Value = -1292.3334
# Row 113
$ws.Range("H113").Value = 86323.75
$ws.Range("J113").Value = 146057.14
$ws.Range("L113").Value = 146057.14
$ws.Range("N113").Value = -152565.14
# Row 122
$ws.Range("H122").Value = 2375.1365
$ws.Range("I122").Value = 2119.1052
$ws.Range("K122").Value = 6357.3156
$ws.Range("M122").Value = -3907.3156
# Row 138
$ws.Range("H138").Value = 3924.122
$ws.Range("I138").Value = 1921
$ws.Range("J138").Value = 4487.5
$ws.Range("K138").Value = 5763
$ws.Range("L138").Value = 13462.5
$ws.Range("M138").Value = -623
$ws.Range("N138").Value = -23742.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 6
$ws.Range("H6").Value = 844333.3
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
# Row 37
$ws.Range("H37").Value = 38996
$ws.Range("I37").Value = 8500
$ws.Range("K37").Value = 8500
$ws.Range("M37").Value = -8227
# Row 45
$ws.Range("H45").Value = 7440.8
$ws.Range("I45").Value = 7156.4443
$ws.Range("J45").Value = 10000
$ws.Range("K45").Value = 7156.4443
$ws.Range("L45").Value = 10000
$ws.Range("M45").Value = -6779.4443
$ws.Range("N45").Value = -10754
# Row 102
$ws.Range("H102").Value = 4267.6816
$ws.Range("I102").Value = 3744.5
$ws.Range("K102").Value = 3744.5
$ws.Range("M102").Value = -2122.5
# Row 125
$ws.Range("H125").Value = 99998
$ws.Range("J125").Value = 99998
$ws.Range("L125").Value = 99998
$ws.Range("N125").Value = -109838
# Row 132
$ws.Range("H132").Value = 5740.32
$ws.Range("I132").Value = 3595.375
$ws.Range("J132").Value = 9553.556
$ws.Range("K132").Value = 10786.125
$ws.Range("L132").Value = 28660.668
$ws.Range("M132").Value = -8256.125
$ws.Range("N132").Value = -33720.66800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 1142.4688
$ws.Range("I94").Value = 1161.1111
$ws.Range("J94").Value = 1041.8
$ws.Range("K94").Value = 1161.1111
$ws.Range("L94").Value = 1041.8
$ws.Range("M94").Value = -710.1111000000001
$ws.Range("N94").Value = -1943.8
# Row 99
$ws.Range("H99").Value = 3657.5625
$ws.Range("I99").Value = 3270.2307
$ws.Range("J99").Value = 5336
$ws.Range("K99").Value = 3270.2307
$ws.Range("L99").Value = 5336
$ws.Range("M99").Value = -1772.2307
$ws.Range("N99").Value = -8332
# Row 135
$ws.Range("H135").Value = 97273.5
$ws.Range("J135").Value = 97273.5
$ws.Range("L135").Value = 97273.5
$ws.Range("N135").Value = -107413.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 7333
$ws.Range("I16").Value = 6666
$ws.Range("J16").Value = 8000
$ws.Range("K16").Value = 6666
$ws.Range("L16").Value = 8000
$ws.Range("M16").Value = -6379
$ws.Range("N16").Value = -8574
# Row 58
$ws.Range("H58").Value = 3720.9302
$ws.Range("I58").Value = 2356.318
$ws.Range("J58").Value = 5150.524
$ws.Range("K58").Value = 2356.318
$ws.Range("L58").Value = 5150.524
$ws.Range("M58").Value = -2153.318
$ws.Range("N58").Value = -5556.524
# Row 107
$ws.Range("H107").Value = 1685.9615
$ws.Range("I107").Value = 1447.6666
$ws.Range("K107").Value = 1447.6666
$ws.Range("M107").Value = 472.3334
# Row 113
$ws.Range("H113").Value = 7333
$ws.Range("I113").Value = 6666
$ws.Range("J113").Value = 8000
$ws.Range("K113").Value = 6666
$ws.Range("L113").Value = 8000
$ws.Range("M113").Value = -4496
$ws.Range("N113").Value = -12340
# Row 136
$ws.Range("H136").Value = 3720.9302
$ws.Range("I136").Value = 2356.318
$ws.Range("J136").Value = 5150.524
$ws.Range("K136").Value = 7068.954000000001
$ws.Range("L136").Value = 15451.572
$ws.Range("M136").Value = -4518.954000000001
$ws.Range("N136").Value = -20551.572

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 1343.5333
$ws.Range("I92").Value = 4419.5
$ws.Range("J92").Value = 870.3077
$ws.Range("K92").Value = 13258.5
$ws.Range("L92").Value = 2610.9231
$ws.Range("M92").Value = -12010.5
$ws.Range("N92").Value = -5106.9231
# Row 140
$ws.Range("H140").Value = 55556660

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 101360
$ws.Range("I7").Value = 101360
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 101360
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -101248
$ws.Range("N7").ClearContents()
# Row 22
$ws.Range("H22").Value = 2402.9795
$ws.Range("I22").Value = 1629.6428
$ws.Range("J22").Value = 3434.0952
$ws.Range("K22").Value = 1629.6428
$ws.Range("L22").Value = 3434.0952
$ws.Range("M22").Value = -1334.6428
$ws.Range("N22").Value = -4024.0952
# Row 27
$ws.Range("H27").Value = 2402.9795
$ws.Range("I27").Value = 1629.6428
$ws.Range("J27").Value = 3434.0952
$ws.Range("K27").Value = 1629.6428
$ws.Range("L27").Value = 3434.0952
$ws.Range("M27").Value = -1522.6428
$ws.Range("N27").Value = -3648.0952
# Row 68
$ws.Range("H68").Value = 1321.1428
$ws.Range("I68").Value = 1321.1428
$ws.Range("K68").Value = 1321.1428
$ws.Range("M68").Value = -572.1428000000001
# Row 71
$ws.Range("H71").Value = 1321.1428
$ws.Range("I71").Value = 1321.1428
$ws.Range("K71").Value = 6605.714
$ws.Range("M71").Value = -2861.714
# Row 94
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
# Row 122
$ws.Range("H122").Value = 50007440
$ws.Range("I122").Value = 55562156
$ws.Range("K122").Value = 166686468
$ws.Range("M122").Value = -166684018
# Row 126
$ws.Range("H126").Value = 101360
$ws.Range("I126").Value = 101360
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 304080
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -301610
$ws.Range("N126").ClearContents()
# Row 127
$ws.Range("H127").Value = 60714.5
$ws.Range("J127").Value = 60714.5
$ws.Range("L127").Value = 60714.5
$ws.Range("N127").Value = -70634.5
# Row 132
$ws.Range("H132").Value = 6371.1055
$ws.Range("I132").Value = 5817.96
$ws.Range("K132").Value = 17453.88
$ws.Range("M132").Value = -14923.88
# Row 133
$ws.Range("H133").Value = 72659.664
$ws.Range("J133").Value = 72659.664
$ws.Range("L133").Value = 72659.664
$ws.Range("N133").Value = -77719.664
# Row 136
$ws.Range("H136").Value = 4256.8667
$ws.Range("I136").Value = 2906.1333
$ws.Range("K136").Value = 8718.3999
$ws.Range("M136").Value = -6168.3999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 28
$ws.Range("H28").Value = 43919.832
$ws.Range("J28").Value = 46703.8
$ws.Range("L28").Value = 46703.8
$ws.Range("N28").Value = -47399.8
# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
# Row 107
$ws.Range("H107").Value = 1827.0714
$ws.Range("I107").Value = 1135.7142
$ws.Range("K107").Value = 3407.1426
$ws.Range("M107").Value = -1487.1426
